# Applies the per-cell value corrections captured in the commit's OOXML diff
# across all eight Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 401.18182
$ws.Range("I2").Value = 201.22223
$ws.Range("J2").Value = 1301
$ws.Range("K2").Value = 201.22223
$ws.Range("L2").Value = 1301
$ws.Range("M2").Value = -88.22223
$ws.Range("N2").Value = -1527

$ws.Range("H132").Value = 3072
$ws.Range("I132").Value = 3055.9546
$ws.Range("J132").Value = 3189.6667
$ws.Range("K132").Value = 9167.863799999999
$ws.Range("L132").Value = 9569.000100000001
$ws.Range("M132").Value = -6637.863799999999
$ws.Range("N132").Value = -14629.0001

$ws.Range("H137").Value = 2052.5334
$ws.Range("I137").Value = 2119.818
$ws.Range("K137").Value = 6359.454000000001
$ws.Range("M137").Value = -3809.454000000001

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 269.16666
$ws.Range("I5").Value = 275
$ws.Range("K5").Value = 275
$ws.Range("M5").Value = -163

$ws.Range("H45").Value = 2404
$ws.Range("I45").Value = 2630.125
$ws.Range("J45").Value = 1499.5
$ws.Range("K45").Value = 2630.125
$ws.Range("L45").Value = 1499.5
$ws.Range("M45").Value = -2253.125
$ws.Range("N45").Value = -2253.5

$ws.Range("H61").Value = 4162.3687
$ws.Range("I61").Value = 3148.7
$ws.Range("J61").Value = 7963.625
$ws.Range("K61").Value = 3148.7
$ws.Range("L61").Value = 7963.625
$ws.Range("M61").Value = -2936.7
$ws.Range("N61").Value = -8387.625

$ws.Range("H62").Value = 32000
$ws.Range("I62").Value = 32000
$ws.Range("K62").Value = 32000
$ws.Range("M62").Value = -31376

$ws.Range("H65").Value = 32000
$ws.Range("I65").Value = 32000
$ws.Range("K65").Value = 96000
$ws.Range("M65").Value = -92880

$ws.Range("H74").Value = 1748
$ws.Range("I74").Value = 1748
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1748
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -874
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1748
$ws.Range("I77").Value = 1748
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 8740
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -4372
$ws.Range("N77").ClearContents()

$ws.Range("H136").Value = 4162.3687
$ws.Range("I136").Value = 3148.7
$ws.Range("J136").Value = 7963.625
$ws.Range("K136").Value = 9446.099999999999
$ws.Range("L136").Value = 23890.875
$ws.Range("M136").Value = -6896.099999999999
$ws.Range("N136").Value = -28990.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 269.16666
$ws.Range("I4").Value = 275
$ws.Range("K4").Value = 275
$ws.Range("M4").Value = -160

$ws.Range("H35").Value = 61428.6
$ws.Range("J35").Value = 63381
$ws.Range("L35").Value = 63381
$ws.Range("N35").Value = -64001

$ws.Range("H105").Value = 3773
$ws.Range("I105").Value = 3869.625
$ws.Range("K105").Value = 3869.625
$ws.Range("M105").Value = -2122.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2638.7144
$ws.Range("I31").Value = 2749.5
$ws.Range("J31").Value = 2594.4
$ws.Range("K31").Value = 2749.5
$ws.Range("L31").Value = 2594.4
$ws.Range("M31").Value = -2454.5
$ws.Range("N31").Value = -3184.4

$ws.Range("H34").Value = 2638.7144
$ws.Range("I34").Value = 2749.5
$ws.Range("J34").Value = 2594.4
$ws.Range("K34").Value = 2749.5
$ws.Range("L34").Value = 2594.4
$ws.Range("M34").Value = -2547.5
$ws.Range("N34").Value = -2998.4

$ws.Range("H63").Value = 81560.234
$ws.Range("J63").Value = 81560.234
$ws.Range("L63").Value = 81560.234
$ws.Range("N63").Value = -82932.234

$ws.Range("H66").Value = 81560.234
$ws.Range("J66").Value = 81560.234
$ws.Range("L66").Value = 244680.702
$ws.Range("N66").Value = -251544.702

$ws.Range("H94").Value = 1916.6666
$ws.Range("I94").Value = 1916.6666
$ws.Range("K94").Value = 1916.6666
$ws.Range("M94").Value = -1465.6666

$ws.Range("H99").Value = 3086.7144
$ws.Range("I99").Value = 3086.7144
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3086.7144
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1588.7144
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 3086.7144
$ws.Range("I126").Value = 3086.7144
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9260.143199999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6790.143199999999
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4959
$ws.Range("I3").Value = 1576.3334
$ws.Range("K3").Value = 4729.0002
$ws.Range("M3").Value = -4617.0002

$ws.Range("H8").Value = 5000
$ws.Range("I8").Value = 5000
$ws.Range("K8").Value = 15000
$ws.Range("M8").Value = -14861

$ws.Range("H39").Value = 6733.1665
$ws.Range("J39").Value = 7318
$ws.Range("L39").Value = 21954
$ws.Range("N39").Value = -22542

$ws.Range("H59").Value = 2499
$ws.Range("I59").Value = 2499
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 7497
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -6957
$ws.Range("N59").ClearContents()

$ws.Range("H75").Value = 585.8
$ws.Range("J75").Value = 699.75
$ws.Range("L75").Value = 2099.25
$ws.Range("N75").Value = -4095.25

$ws.Range("H78").Value = 585.8
$ws.Range("J78").Value = 699.75
$ws.Range("L78").Value = 6297.75
$ws.Range("N78").Value = -16281.75

$ws.Range("H92").Value = 210.1
$ws.Range("I92").Value = 221.2
$ws.Range("J92").Value = 199
$ws.Range("K92").Value = 663.5999999999999
$ws.Range("L92").Value = 597
$ws.Range("M92").Value = 584.4000000000001
$ws.Range("N92").Value = -3093

$ws.Range("H131").Value = 572.6667
$ws.Range("I131").Value = 601.75
$ws.Range("J131").Value = 340
$ws.Range("K131").Value = 1805.25
$ws.Range("L131").Value = 1020
$ws.Range("M131").Value = 3234.75
$ws.Range("N131").Value = -11100

$ws.Range("H136").Value = 7806.75
$ws.Range("I136").Value = 6345.857
$ws.Range("K136").Value = 19037.571
$ws.Range("M136").Value = -13937.571

$ws.Range("H138").Value = 6697.2144
$ws.Range("I138").Value = 6697.2144
$ws.Range("K138").Value = 20091.6432
$ws.Range("M138").Value = -14951.6432

$ws.Range("H139").Value = 4893.1177
$ws.Range("I139").Value = 4812.2
$ws.Range("K139").Value = 14436.6
$ws.Range("M139").Value = -9296.599999999999

$ws.Range("H141").Value = 7818.3335
$ws.Range("I141").Value = 7818.3335
$ws.Range("K141").Value = 23455.0005
$ws.Range("M141").Value = -18275.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 14305.5
$ws.Range("J92").Value = 14305.5
$ws.Range("L92").Value = 14305.5
$ws.Range("N92").Value = -18049.5

$ws.Range("H132").Value = 2677.4285
$ws.Range("I132").Value = 2575.6924
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7727.0772
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5197.0772
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 730.3
$ws.Range("I16").Value = 614.2857
$ws.Range("J16").Value = 1001
$ws.Range("K16").Value = 614.2857
$ws.Range("L16").Value = 1001
$ws.Range("M16").Value = -444.2857
$ws.Range("N16").Value = -1341

$ws.Range("H22").Value = 599.5
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 599.5
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 599.5
$ws.Range("N22").Value = -1189.5
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 599.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 599.5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 599.5
$ws.Range("N27").Value = -813.5
$ws.Range("M27").ClearContents()

$ws.Range("H40").Value = 2025.2632
$ws.Range("I40").Value = 1911.6875
$ws.Range("K40").Value = 1911.6875
$ws.Range("M40").Value = -1775.6875

$ws.Range("H46").Value = 1090.2727
$ws.Range("J46").Value = 1279
$ws.Range("L46").Value = 1279
$ws.Range("N46").Value = -1655

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 6414.231
$ws.Range("I136").Value = 5398.636
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 16195.908
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -13645.908
$ws.Range("N136").Value = -41100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4055.2856
$ws.Range("I96").Value = 3374.5
$ws.Range("J96").Value = 4963
$ws.Range("K96").Value = 3374.5
$ws.Range("L96").Value = 4963
$ws.Range("M96").Value = -2001.5
$ws.Range("N96").Value = -7709

$ws.Range("H113").Value = 717.7273
$ws.Range("I113").Value = 539.5
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 7500
$ws.Range("M113").Value = 551.5
$ws.Range("N113").Value = -11840

Write-Host "Applied all cell updates"
